$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 90
$ws.Range("D90").Value = 44664
$ws.Range("K90").Value = "Sin especificar"
$ws.Range("L90").Value = "3a plateado"
$ws.Range("M90").Value = 200
$ws.Range("N90").Value = 24000
$ws.Range("O90").Value = 25000
$ws.Range("P90").Value = 24500
$ws.Range("Q90").Value = "$/caja 20 kilos"
$ws.Range("R90").Value = "Región de Coquimbo"
$ws.Range("S90").Value = 1225
$ws.Range("T90").Value = 20

# Row 91
$ws.Range("D91").Value = 44175
$ws.Range("K91").Value = "Sin especificar"
$ws.Range("L91").Value = "3a amarillo"
$ws.Range("M91").Value = 250
$ws.Range("N91").Value = 13000
$ws.Range("O91").Value = 14000
$ws.Range("P91").Value = 13500
$ws.Range("Q91").Value = "$/caja 20 kilos"
$ws.Range("R91").Value = "Región Metropolitana"
$ws.Range("S91").Value = 675
$ws.Range("T91").Value = 20

# Row 92
$ws.Range("D92").Value = 44169
$ws.Range("K92").Value = "Sutil De Gase"
$ws.Range("L92").Value = "Primera"
$ws.Range("M92").Value = 270
$ws.Range("N92").Value = 29000
$ws.Range("O92").Value = 31000
$ws.Range("P92").Value = 30000
$ws.Range("Q92").Value = "$/caja 24 kilos"
$ws.Range("R92").Value = "Perú"
$ws.Range("S92").Value = 1250
$ws.Range("T92").Value = 24

# Row 93
$ws.Range("D93").Value = 44169
$ws.Range("K93").Value = "Tahití"
$ws.Range("L93").Value = "Primera"
$ws.Range("M93").Value = 360
$ws.Range("N93").Value = 23000
$ws.Range("O93").Value = 24000
$ws.Range("P93").Value = 23500
$ws.Range("Q93").Value = "$/caja 24 kilos"
$ws.Range("R93").Value = "Perú"
$ws.Range("S93").Value = 979
$ws.Range("T93").Value = 24

# Row 94
$ws.Range("D94").Value = 44643
$ws.Range("K94").Value = "Sin especificar"
$ws.Range("L94").Value = "2a plateado"
$ws.Range("M94").Value = 200
$ws.Range("N94").Value = 30000
$ws.Range("O94").Value = 32000
$ws.Range("P94").Value = 31000
$ws.Range("Q94").Value = "$/caja 20 kilos"
$ws.Range("R94").Value = "Región de Coquimbo"
$ws.Range("S94").Value = 1550
$ws.Range("T94").Value = 20

# Row 95
$ws.Range("D95").Value = 44643
$ws.Range("K95").Value = "Sutil De Gase"
$ws.Range("L95").Value = "Primera"
$ws.Range("M95").Value = 270
$ws.Range("N95").Value = 33000
$ws.Range("O95").Value = 34000
$ws.Range("P95").Value = 33500
$ws.Range("Q95").Value = "$/caja 24 kilos"
$ws.Range("R95").Value = "Región de Arica y Parinacota"
$ws.Range("S95").Value = 1396
$ws.Range("T95").Value = 24

# Row 96
$ws.Range("D96").Value = 44643
$ws.Range("K96").Value = "Tahití"
$ws.Range("L96").Value = "Primera"
$ws.Range("M96").Value = 300
$ws.Range("N96").Value = 36000
$ws.Range("O96").Value = 37000
$ws.Range("P96").Value = 36500
$ws.Range("Q96").Value = "$/caja 24 kilos"
$ws.Range("R96").Value = "Región de Arica y Parinacota"
$ws.Range("S96").Value = 1521
$ws.Range("T96").Value = 24

# Row 97
$ws.Range("D97").Value = 44319
$ws.Range("K97").Value = "Sutil De Gase"
$ws.Range("L97").Value = "Primera"
$ws.Range("M97").Value = 250
$ws.Range("N97").Value = 27000
$ws.Range("O97").Value = 28000
$ws.Range("P97").Value = 27500
$ws.Range("Q97").Value = "$/caja 24 kilos"
$ws.Range("R97").Value = "Perú"
$ws.Range("S97").Value = 1146
$ws.Range("T97").Value = 24

# Row 98
$ws.Range("D98").Value = 44319
$ws.Range("K98").Value = "Tahití"
$ws.Range("L98").Value = "Primera"
$ws.Range("M98").Value = 300
$ws.Range("N98").Value = 23000
$ws.Range("O98").Value = 24000
$ws.Range("P98").Value = 23500
$ws.Range("Q98").Value = "$/caja 24 kilos"
$ws.Range("R98").Value = "Perú"
$ws.Range("S98").Value = 979
$ws.Range("T98").Value = 24

# Row 99
$ws.Range("D99").Value = 44300
$ws.Range("K99").Value = "Sin especificar"
$ws.Range("L99").Value = "2a amarillo"
$ws.Range("M99").Value = 270
$ws.Range("N99").Value = 20000
$ws.Range("O99").Value = 21000
$ws.Range("P99").Value = 20500
$ws.Range("Q99").Value = "$/caja 20 kilos"
$ws.Range("R99").Value = "Región de Coquimbo"
$ws.Range("S99").Value = 1025
$ws.Range("T99").Value = 20

# Row 100
$ws.Range("D100").Value = 44624
$ws.Range("K100").Value = "Sutil De Gase"
$ws.Range("L100").Value = "Primera"
$ws.Range("M100").Value = 200
$ws.Range("N100").Value = 46000
$ws.Range("O100").Value = 47000
$ws.Range("P100").Value = 46500
$ws.Range("Q100").Value = "$/caja 24 kilos"
$ws.Range("R100").Value = "Perú"
$ws.Range("S100").Value = 1938
$ws.Range("T100").Value = 24

# Row 101
$ws.Range("D101").Value = 44624
$ws.Range("K101").Value = "Tahití"
$ws.Range("L101").Value = "Primera"
$ws.Range("M101").Value = 300
$ws.Range("N101").Value = 45000
$ws.Range("O101").Value = 46000
$ws.Range("P101").Value = 45500
$ws.Range("Q101").Value = "$/caja 24 kilos"
$ws.Range("R101").Value = "Perú"
$ws.Range("S101").Value = 1896
$ws.Range("T101").Value = 24

# Row 102
$ws.Range("D102").Value = 44239
$ws.Range("K102").Value = "Sutil De Gase"
$ws.Range("L102").Value = "Primera"
$ws.Range("M102").Value = 250
$ws.Range("N102").Value = 22000
$ws.Range("O102").Value = 23000
$ws.Range("P102").Value = 22500
$ws.Range("Q102").Value = "$/caja 24 kilos"
$ws.Range("R102").Value = "Perú"
$ws.Range("S102").Value = 938
$ws.Range("T102").Value = 24

# Row 103
$ws.Range("D103").Value = 44239
$ws.Range("K103").Value = "Tahití"
$ws.Range("L103").Value = "Primera"
$ws.Range("M103").Value = 400
$ws.Range("N103").Value = 20000
$ws.Range("O103").Value = 21000
$ws.Range("P103").Value = 20500
$ws.Range("Q103").Value = "$/caja 24 kilos"
$ws.Range("R103").Value = "Perú"
$ws.Range("S103").Value = 854
$ws.Range("T103").Value = 24

# Row 104
$ws.Range("D104").Value = 44329
$ws.Range("K104").Value = "Sin especificar"
$ws.Range("L104").Value = "2a amarillo"
$ws.Range("M104").Value = 250
$ws.Range("N104").Value = 25000
$ws.Range("O104").Value = 26000
$ws.Range("P104").Value = 25500
$ws.Range("Q104").Value = "$/caja 20 kilos"
$ws.Range("R104").Value = "Región de Coquimbo"
$ws.Range("S104").Value = 1275
$ws.Range("T104").Value = 20

# Row 105
$ws.Range("D105").Value = 44630
$ws.Range("K105").Value = "Sin especificar"
$ws.Range("L105").Value = "3a amarillo"
$ws.Range("M105").Value = 250
$ws.Range("N105").Value = 28000
$ws.Range("O105").Value = 30000
$ws.Range("P105").Value = 29000
$ws.Range("Q105").Value = "$/caja 20 kilos"
$ws.Range("R105").Value = "Región de Coquimbo"
$ws.Range("S105").Value = 1450
$ws.Range("T105").Value = 20

# Row 106
$ws.Range("D106").Value = 44421
$ws.Range("K106").Value = "Sutil De Gase"
$ws.Range("L106").Value = "Primera"
$ws.Range("M106").Value = 250
$ws.Range("N106").Value = 31000
$ws.Range("O106").Value = 32000
$ws.Range("P106").Value = 31500
$ws.Range("Q106").Value = "$/caja 24 kilos"
$ws.Range("R106").Value = "Perú"
$ws.Range("S106").Value = 1312
$ws.Range("T106").Value = 24

# Row 107
$ws.Range("D107").Value = 44421
$ws.Range("K107").Value = "Tahití"
$ws.Range("L107").Value = "Primera"
$ws.Range("M107").Value = 300
$ws.Range("N107").Value = 30000
$ws.Range("O107").Value = 31000
$ws.Range("P107").Value = 30500
$ws.Range("Q107").Value = "$/caja 24 kilos"
$ws.Range("R107").Value = "Perú"
$ws.Range("S107").Value = 1271
$ws.Range("T107").Value = 24

# Row 108
$ws.Range("D108").Value = 44298
$ws.Range("K108").Value = "Sutil De Gase"
$ws.Range("L108").Value = "Primera"
$ws.Range("M108").Value = 160
$ws.Range("N108").Value = 31000
$ws.Range("O108").Value = 32000
$ws.Range("P108").Value = 31500
$ws.Range("Q108").Value = "$/caja 24 kilos"
$ws.Range("R108").Value = "Perú"
$ws.Range("S108").Value = 1312
$ws.Range("T108").Value = 24

# Row 109
$ws.Range("D109").Value = 44298
$ws.Range("K109").Value = "Tahití"
$ws.Range("L109").Value = "Primera"
$ws.Range("M109").Value = 300
$ws.Range("N109").Value = 27000
$ws.Range("O109").Value = 28000
$ws.Range("P109").Value = 27500
$ws.Range("Q109").Value = "$/caja 24 kilos"
$ws.Range("R109").Value = "Perú"
$ws.Range("S109").Value = 1146
$ws.Range("T109").Value = 24

# Row 110
$ws.Range("D110").Value = 44469
$ws.Range("K110").Value = "Sin especificar"
$ws.Range("L110").Value = "2a amarillo"
$ws.Range("M110").Value = 250
$ws.Range("N110").Value = 10000
$ws.Range("O110").Value = 11000
$ws.Range("P110").Value = 10500
$ws.Range("Q110").Value = "$/caja 20 kilos"
$ws.Range("R110").Value = "Región Metropolitana"
$ws.Range("S110").Value = 525
$ws.Range("T110").Value = 20

# Row 111
$ws.Range("D111").Value = 44267
$ws.Range("K111").Value = "Sutil De Gase"
$ws.Range("L111").Value = "Primera"
$ws.Range("M111").Value = 250
$ws.Range("N111").Value = 26000
$ws.Range("O111").Value = 27000
$ws.Range("P111").Value = 26500
$ws.Range("Q111").Value = "$/caja 24 kilos"
$ws.Range("R111").Value = "Perú"
$ws.Range("S111").Value = 1104
$ws.Range("T111").Value = 24

# Row 112
$ws.Range("D112").Value = 44267
$ws.Range("K112").Value = "Tahití"
$ws.Range("L112").Value = "Primera"
$ws.Range("M112").Value = 300
$ws.Range("N112").Value = 23000
$ws.Range("O112").Value = 24000
$ws.Range("P112").Value = 23500
$ws.Range("Q112").Value = "$/caja 24 kilos"
$ws.Range("R112").Value = "Perú"
$ws.Range("S112").Value = 979
$ws.Range("T112").Value = 24

# Row 113
$ws.Range("D113").Value = 44475
$ws.Range("K113").Value = "Sin especificar"
$ws.Range("L113").Value = "1a amarillo"
$ws.Range("M113").Value = 300
$ws.Range("N113").Value = 13000
$ws.Range("O113").Value = 14000
$ws.Range("P113").Value = 13500
$ws.Range("Q113").Value = "$/caja 20 kilos"
$ws.Range("R113").Value = "Región de Coquimbo"
$ws.Range("S113").Value = 675
$ws.Range("T113").Value = 20

# Row 114
$ws.Range("D114").Value = 44333
$ws.Range("K114").Value = "Tahití"
$ws.Range("L114").Value = "Primera"
$ws.Range("M114").Value = 250
$ws.Range("N114").Value = 27000
$ws.Range("O114").Value = 28000
$ws.Range("P114").Value = 27500
$ws.Range("Q114").Value = "$/caja 24 kilos"
$ws.Range("R114").Value = "Perú"
$ws.Range("S114").Value = 1146
$ws.Range("T114").Value = 24

# Row 115
$ws.Range("D115").Value = 44218
$ws.Range("K115").Value = "Sutil De Gase"
$ws.Range("L115").Value = "Primera"
$ws.Range("M115").Value = 300
$ws.Range("N115").Value = 21000
$ws.Range("O115").Value = 22000
$ws.Range("P115").Value = 21500
$ws.Range("Q115").Value = "$/caja 24 kilos"
$ws.Range("R115").Value = "Perú"
$ws.Range("S115").Value = 896
$ws.Range("T115").Value = 24

# Row 116
$ws.Range("D116").Value = 44218
$ws.Range("K116").Value = "Tahití"
$ws.Range("L116").Value = "Primera"
$ws.Range("M116").Value = 300
$ws.Range("N116").Value = 21000
$ws.Range("O116").Value = 22000
$ws.Range("P116").Value = 21500
$ws.Range("Q116").Value = "$/caja 24 kilos"
$ws.Range("R116").Value = "Perú"
$ws.Range("S116").Value = 896
$ws.Range("T116").Value = 24

# Row 117
$ws.Range("D117").Value = 44340
$ws.Range("K117").Value = "Tahití"
$ws.Range("L117").Value = "Primera"
$ws.Range("M117").Value = 250
$ws.Range("N117").Value = 25000
$ws.Range("O117").Value = 26000
$ws.Range("P117").Value = 25500
$ws.Range("Q117").Value = "$/caja 24 kilos"
$ws.Range("R117").Value = "Perú"
$ws.Range("S117").Value = 1062
$ws.Range("T117").Value = 24

